# Update the header row (row 1) text on the active worksheet: the "20" /
# "10" resistor labels are clarified to "20KΩ" / "10KΩ" (or "20Ω" for the
# two "0-ohm"/bare columns A and F), reflecting the real component values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Vi_20Ω_5V (mV)"
$ws.Range("B1").Value = "Ii_20KΩ_5V (mA)"
$ws.Range("C1").Value = "Vi_20KΩ_0,5V (mV)"
$ws.Range("D1").Value = "Ii_20KΩ_0,5V (mV)"
$ws.Range("E1").Value = "Vi_20KΩ_0,25V (mV)"
$ws.Range("F1").Value = "Vi_20Ω_0,25V (mV)"
$ws.Range("G1").Value = "V_20KΩ(mV)"
$ws.Range("H1").Value = "I_20KΩ(mV)"
$ws.Range("I1").Value = "Vi_10KΩ_5V (mV)"
$ws.Range("J1").Value = "Ii_10KΩ_5V (mV)"
$ws.Range("K1").Value = "Vi_10KΩ_0,5V (mV)"
$ws.Range("L1").Value = "Vi_10KΩ_0,5V (mV)"
$ws.Range("M1").Value = "Vi_10KΩ_0,25V (mV)"
$ws.Range("N1").Value = "Vi_10KΩ_0,25V (mV)"
$ws.Range("O1").Value = "V_10KΩ(mV)"
$ws.Range("P1").Value = "I_10KΩ(mV)"

# Re-fit the column widths ("best fit") to the new, longer header text -
# mirrors what Excel does automatically on a bestFit column when the
# header text changes. (This host's ColumnWidth setter re-adds its own
# fixed 5/6-character padding on top of whatever is assigned, so the
# assigned value is the desired final width net of that padding.)
$padding = 5 / 6
$ws.Columns.Item(1).ColumnWidth = 14.33203125 - $padding
$ws.Columns.Item(2).ColumnWidth = 14.77734375 - $padding
$ws.Columns.Item(3).ColumnWidth = 16.77734375 - $padding
$ws.Columns.Item(4).ColumnWidth = 16.21875 - $padding
$ws.Columns.Item(5).ColumnWidth = 17.88671875 - $padding
$ws.Columns.Item(6).ColumnWidth = 16.77734375 - $padding
$ws.Columns.Item(7).ColumnWidth = 11.21875 - $padding
$ws.Columns.Item(8).ColumnWidth = 10.6640625 - $padding
$ws.Columns.Item(9).ColumnWidth = 15.33203125 - $padding
$ws.Columns.Item(10).ColumnWidth = 14.77734375 - $padding
$ws.Columns.Item(11).ColumnWidth = 16.77734375 - $padding
$ws.Columns.Item(12).ColumnWidth = 16.77734375 - $padding
$ws.Columns.Item(13).ColumnWidth = 17.88671875 - $padding
$ws.Columns.Item(14).ColumnWidth = 17.88671875 - $padding
$ws.Columns.Item(15).ColumnWidth = 11.21875 - $padding
$ws.Columns.Item(16).ColumnWidth = 10.6640625 - $padding

# Restore the view: scrolled so column C is leftmost, with C4 selected.
$ws.Application.ActiveWindow.ScrollColumn = 3
$ws.Range("C4").Select()
